$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E23").Value = 120
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = "体力存储上限"
$ws.Range("E24").Value = 10000
$ws.Range("K23").Select() | Out-Null
